$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two bugs fixed: both hyperlink cells C7 and C8 shared the placeholder text "t".
# Give each its own proper label.
$ws.Range("C7").Value = "link"
$ws.Range("C8").Value = "link2"
